$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# Add the three new adjustment columns (ADJ1, ADJ2, ADJ3) to the header row.
$ws.Range("S1").Value = "ADJ1"
$ws.Range("T1").Value = "ADJ2"
$ws.Range("U1").Value = "ADJ3"

# Move the active selection to the newly added header cells, matching the
# state captured after the edit was made in Excel.
$ws.Range("S1:U1").Select()
